$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '26.385.17'
$ws.Range('E2').Value = '  -0.41%  '

$ws.Range('D3').Value = '1.724.09'
$ws.Range('E3').Value = '  -0.44%  '

Set-TextValue $ws.Range('D4') '0.9990'
$ws.Range('E4').Value = '  -0.11%  '

Set-TextValue $ws.Range('D5') '242.53'
$ws.Range('E5').Value = '  -1.92%  '

Set-TextValue $ws.Range('D6') '0.9995'
$ws.Range('E6').Value = '  -0.07%  '

Set-TextValue $ws.Range('D7') '0.4885'
$ws.Range('E7').Value = '  +0.47%  '

$ws.Range('E8').Value = '  -3.13%  '

Set-TextValue $ws.Range('D9') '0.06188'

$ws.Range('D10').Value = '1.726.50'
$ws.Range('E10').Value = '  -0.31%  '

Set-TextValue $ws.Range('D11') '0.06979'
$ws.Range('E11').Value = '  -1.26%  '

$ws.Range('E12').Value = '  -1.02%  '

Set-TextValue $ws.Range('D13') '4.520'
$ws.Range('E13').Value = '  -2.57%  '

Set-TextValue $ws.Range('D14') '0.5977'
$ws.Range('E14').Value = '  -2.00%  '

$ws.Range('E15').Value = '  -0.47%  '

Set-TextValue $ws.Range('D16') '0.9991'
$ws.Range('E16').Value = '  -0.10%  '

$ws.Range('D17').Value = '26.384.90'
$ws.Range('E17').Value = '  -0.41%  '

Set-TextValue $ws.Range('D18') '0.9993'
$ws.Range('E18').Value = '  -0.12%  '

Set-TextValue $ws.Range('D19') '0.000007177'
$ws.Range('E19').Value = '  +0.64%  '

$ws.Range('D21').Value = '1.939.97'
$ws.Range('E21').Value = '  -0.98%  '

Set-TextValue $ws.Range('D22') '4.450'
$ws.Range('E22').Value = '  -1.66%  '

Set-TextValue $ws.Range('D23') '8.487'
$ws.Range('E23').Value = '  -3.29%  '

Set-TextValue $ws.Range('D24') '5.105'
$ws.Range('E24').Value = '  -2.92%  '

Set-TextValue $ws.Range('D25') '138.10'
$ws.Range('E25').Value = '  -0.46%  '

$ws.Range('E26').Value = '  -1.29%  '

Set-TextValue $ws.Range('D27') '1.399'
$ws.Range('E27').Value = '  -0.33%  '

Set-TextValue $ws.Range('D28') '106.33'
$ws.Range('E28').Value = '  -1.76%  '

Set-TextValue $ws.Range('D29') '1.728'
$ws.Range('E29').Value = '  -2.93%  '

Set-TextValue $ws.Range('D30') '3.906'
$ws.Range('E30').Value = '  -1.81%  '

Set-TextValue $ws.Range('D31') '0.08010'
$ws.Range('E31').Value = '  -0.15%  '

Set-TextValue $ws.Range('D32') '3.658'
$ws.Range('E32').Value = '  -0.97%  '

Set-TextValue $ws.Range('D33') '0.04499'
$ws.Range('E33').Value = '  -1.85%  '

$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D34') '2.604'
$ws.Range('E34').Value = '  -0.50%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D35') '0.9973'
$ws.Range('E35').Value = '  -0.71%  '

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D36') '0.6228'
$ws.Range('E36').Value = '  -2.07%  '

$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D37') '0.9304'
$ws.Range('E37').Value = '  +3.41%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D38') '1.964'
$ws.Range('E38').Value = '  -3.07%  '

$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D39') '2.386'
$ws.Range('E39').Value = '  -0.53%  '

$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range('D40') '0.9987'
$ws.Range('E40').Value = '  -0.53%  '

$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D41') '0.01473'
$ws.Range('E41').Value = '  -2.32%  '

$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range('D42') '100.17'
$ws.Range('E42').Value = '  -1.35%  '

$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D43') '5.453'
$ws.Range('E43').Value = '  +0.02%  '

$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D44') '0.3837'
$ws.Range('E44').Value = '  -1.43%  '

$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D45') '6.903'
$ws.Range('E45').Value = '  -0.99%  '

$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D46') '0.1164'
$ws.Range('E46').Value = '  -1.82%  '

$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D47') '0.05367'
$ws.Range('E47').Value = '  -0.32%  '

$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue $ws.Range('D48') '30.13'
$ws.Range('E48').Value = '  -1.63%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D49') '7.678'
$ws.Range('E49').Value = '  -2.09%  '

$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D50') '1.225'
$ws.Range('E50').Value = '  -1.94%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D51') '50.91'
$ws.Range('E51').Value = '  -0.99%  '
